$wb = $excel.ActiveWorkbook

# Sheet ALC, row 33 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 710.75
$ws.Range("I33").Value = 851.6923
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 851.6923
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -622.6923
$ws.Range("N33").Value = -558

# Sheet ALC, row 43 (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 604.9231
$ws.Range("I43").Value = 553.4286
$ws.Range("J43").Value = 665
$ws.Range("K43").Value = 553.4286
$ws.Range("L43").Value = 665
$ws.Range("M43").Value = -484.4286
$ws.Range("N43").Value = -803

# Sheet ALC, row 103 (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1616.6666
$ws.Range("I103").Value = 1687.5
$ws.Range("J103").Value = 1560
$ws.Range("K103").Value = 5062.5
$ws.Range("L103").Value = 4680
$ws.Range("M103").Value = -4476.5
$ws.Range("N103").Value = -5852

# Sheet ALC, row 125 (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 4439.3335
$ws.Range("I125").Value = 800
$ws.Range("J125").Value = 6259
$ws.Range("K125").Value = 7200
$ws.Range("L125").Value = 56331
$ws.Range("M125").Value = -4740
$ws.Range("N125").Value = -61251

# Sheet ARM, row 32 (hunk 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4817.983
$ws.Range("I32").Value = 4977.2183
$ws.Range("J32").Value = 2628.5
$ws.Range("K32").Value = 4977.2183
$ws.Range("L32").Value = 2628.5
$ws.Range("M32").Value = -4690.2183
$ws.Range("N32").Value = -3202.5

# Sheet ARM, row 61 (hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2771.9395
$ws.Range("I61").Value = 1866.5714
$ws.Range("J61").Value = 3439.0527
$ws.Range("K61").Value = 1866.5714
$ws.Range("L61").Value = 3439.0527
$ws.Range("M61").Value = -1654.5714
$ws.Range("N61").Value = -3863.0527

# Sheet ARM, row 102 (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3000
$ws.Range("I102").Value = 3000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1378

# Sheet ARM, row 110 (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2246.238
$ws.Range("I110").Value = 704.55554
$ws.Range("J110").Value = 3402.5
$ws.Range("K110").Value = 704.55554
$ws.Range("L110").Value = 3402.5
$ws.Range("M110").Value = 1340.44446
$ws.Range("N110").Value = -7492.5

# Sheet ARM, row 122 (hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1783
$ws.Range("I122").Value = 1386.359
$ws.Range("J122").Value = 3992.8572
$ws.Range("K122").Value = 4159.076999999999
$ws.Range("L122").Value = 11978.5716
$ws.Range("M122").Value = -1709.076999999999
$ws.Range("N122").Value = -16878.5716

# Sheet ARM, row 136 (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2771.9395
$ws.Range("I136").Value = 1866.5714
$ws.Range("J136").Value = 3439.0527
$ws.Range("K136").Value = 5599.7142
$ws.Range("L136").Value = 10317.1581
$ws.Range("M136").Value = -3049.7142
$ws.Range("N136").Value = -15417.1581

# Sheet BSM, row 80 (hunk 10)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 194.95238
$ws.Range("I80").Value = 201.71428
$ws.Range("J80").Value = 181.42857
$ws.Range("K80").Value = 201.71428
$ws.Range("L80").Value = 181.42857
$ws.Range("M80").Value = 796.28572
$ws.Range("N80").Value = -2177.42857

# Sheet BSM, row 83 (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 194.95238
$ws.Range("I83").Value = 201.71428
$ws.Range("J83").Value = 181.42857
$ws.Range("K83").Value = 1008.5714
$ws.Range("L83").Value = 907.1428500000001
$ws.Range("M83").Value = 3983.4286
$ws.Range("N83").Value = -10891.14285

# Sheet CRP, row 16 (hunk 12)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3910.1538
$ws.Range("I16").Value = 3283.2
$ws.Range("J16").Value = 6000
$ws.Range("K16").Value = 3283.2
$ws.Range("L16").Value = 6000
$ws.Range("M16").Value = -2996.2
$ws.Range("N16").Value = -6574

# Sheet CRP, row 31 (hunk 13)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6669138.5
$ws.Range("I31").Value = 1781.4642
$ws.Range("J31").Value = 15154866
$ws.Range("K31").Value = 1781.4642
$ws.Range("L31").Value = 15154866
$ws.Range("M31").Value = -1486.4642
$ws.Range("N31").Value = -15155456

# Sheet CRP, row 34 (hunk 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6669138.5
$ws.Range("I34").Value = 1781.4642
$ws.Range("J34").Value = 15154866
$ws.Range("K34").Value = 1781.4642
$ws.Range("L34").Value = 15154866
$ws.Range("M34").Value = -1579.4642
$ws.Range("N34").Value = -15155270

# Sheet CRP, row 48 (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 6571.4287
$ws.Range("I48").Value = 3000
$ws.Range("J48").Value = 8000
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 8000
$ws.Range("M48").Value = -2524
$ws.Range("N48").Value = -8952

# Sheet CRP, row 113 (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3910.1538
$ws.Range("I113").Value = 3283.2
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 3283.2
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -1113.2
$ws.Range("N113").Value = -10340

# Sheet CRP, row 132 (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2790.1538
$ws.Range("I132").Value = 1878.8334
$ws.Range("J132").Value = 3571.2856
$ws.Range("K132").Value = 5636.5002
$ws.Range("L132").Value = 10713.8568
$ws.Range("M132").Value = -3106.5002
$ws.Range("N132").Value = -15773.8568

# Sheet CUL, row 113 (hunk 18)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 824.76
$ws.Range("I113").Value = 994.6667
$ws.Range("J113").Value = 701.7241
$ws.Range("K113").Value = 2984.0001
$ws.Range("L113").Value = 2105.1723
$ws.Range("M113").Value = -814.0001000000002
$ws.Range("N113").Value = -6445.1723

# Sheet GSM, row 62 (hunk 19)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 11000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11000
$ws.Range("N62").Value = -12372

# Sheet GSM, row 65 (hunk 20)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 11000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 33000
$ws.Range("N65").Value = -39864

# Sheet GSM, row 69 (hunk 21)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 25000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 25000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498

# Sheet GSM, row 72 (hunk 22)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 25000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 25000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488

# Sheet GSM, row 95 (hunk 23)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 6272
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 6272
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 6272
$ws.Range("N95").Value = -11764

# Sheet GSM, row 113 (hunk 24)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2032.3
$ws.Range("I113").Value = 2230
$ws.Range("J113").Value = 1571
$ws.Range("K113").Value = 2230
$ws.Range("L113").Value = 1571
$ws.Range("M113").Value = -60
$ws.Range("N113").Value = -5911

# Sheet GSM, row 132 (hunk 25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4276.2856
$ws.Range("I132").Value = 5345.3335
$ws.Range("J132").Value = 3474.5
$ws.Range("K132").Value = 16036.0005
$ws.Range("L132").Value = 10423.5
$ws.Range("M132").Value = -13506.0005
$ws.Range("N132").Value = -15483.5

# Sheet LTW, row 22 (hunk 26)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 760.7059
$ws.Range("I22").Value = 345.44446
$ws.Range("J22").Value = 1227.875
$ws.Range("K22").Value = 345.44446
$ws.Range("L22").Value = 1227.875
$ws.Range("M22").Value = -50.44445999999999
$ws.Range("N22").Value = -1817.875

# Sheet LTW, row 27 (hunk 27)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 760.7059
$ws.Range("I27").Value = 345.44446
$ws.Range("J27").Value = 1227.875
$ws.Range("K27").Value = 345.44446
$ws.Range("L27").Value = 1227.875
$ws.Range("M27").Value = -238.44446
$ws.Range("N27").Value = -1441.875

# Sheet LTW, row 46 (hunk 28)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 477184.75
$ws.Range("I46").Value = 806.4286
$ws.Range("J46").Value = 1429941.4
$ws.Range("K46").Value = 806.4286
$ws.Range("L46").Value = 1429941.4
$ws.Range("M46").Value = -618.4286
$ws.Range("N46").Value = -1430317.4

# Sheet LTW, row 55 (hunk 29)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 267.62857
$ws.Range("I55").Value = 219.40909
$ws.Range("J55").Value = 349.23077
$ws.Range("K55").Value = 219.40909
$ws.Range("L55").Value = 349.23077
$ws.Range("M55").Value = -46.40908999999999
$ws.Range("N55").Value = -695.23077

# Sheet LTW, row 61 (hunk 30)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5206.9375
$ws.Range("I61").Value = 4070.1
$ws.Range("J61").Value = 7101.6665
$ws.Range("K61").Value = 4070.1
$ws.Range("L61").Value = 7101.6665
$ws.Range("M61").Value = -3868.1
$ws.Range("N61").Value = -7505.6665

# Sheet LTW, row 100 (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3441.2
$ws.Range("I100").Value = 3150
$ws.Range("J100").Value = 3547.0908
$ws.Range("K100").Value = 3150
$ws.Range("L100").Value = 3547.0908
$ws.Range("M100").Value = -2609
$ws.Range("N100").Value = -4629.0908

# Sheet LTW, row 113 (hunk 32)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5206.9375
$ws.Range("I113").Value = 4070.1
$ws.Range("J113").Value = 7101.6665
$ws.Range("K113").Value = 4070.1
$ws.Range("L113").Value = 7101.6665
$ws.Range("M113").Value = -1900.1
$ws.Range("N113").Value = -11441.6665

# Sheet WVR, row 54 (hunk 33)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11996.333
$ws.Range("I54").Value = 6000
$ws.Range("J54").Value = 14994.5
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 14994.5
$ws.Range("M54").Value = -5480
$ws.Range("N54").Value = -16034.5

# Sheet WVR, row 122 (hunk 34)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5536.268
$ws.Range("I122").Value = 6187.7915
$ws.Range("J122").Value = 4616.4707
$ws.Range("K122").Value = 18563.3745
$ws.Range("L122").Value = 13849.4121
$ws.Range("M122").Value = -16113.3745
$ws.Range("N122").Value = -18749.4121

# Sheet WVR, row 132 (hunk 35)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2437.353
$ws.Range("I132").Value = 1918.35
$ws.Range("J132").Value = 3178.7856
$ws.Range("K132").Value = 5755.049999999999
$ws.Range("L132").Value = 9536.356800000001
$ws.Range("M132").Value = -3225.049999999999
$ws.Range("N132").Value = -14596.3568
